# Update the "想去人数" (F column) figures across the relevant worksheets.
# Values/rows derived from the authoritative diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 349
    4  = 1263
    6  = 27
    9  = 146
    10 = 3497
    11 = 133
    12 = 88
    13 = 69
    14 = 43
    15 = 55
    16 = 601
    17 = 96
    18 = 748
    19 = 210
    22 = 59
    24 = 2660
    25 = 5160
    28 = 478
    29 = 3073
    31 = 2252
    35 = 121
    36 = 179
    38 = 16
    40 = 801
    42 = 2
    44 = 38
    45 = 484
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 73

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 349
    4  = 1263
    6  = 27
    9  = 146
    10 = 3497
    11 = 133
    12 = 88
    13 = 69
    14 = 73
    15 = 43
    16 = 55
    17 = 601
    18 = 96
    19 = 748
    20 = 210
    23 = 59
    25 = 2660
    26 = 5160
    29 = 478
    30 = 3073
    32 = 2252
    36 = 121
    37 = 179
    39 = 16
    41 = 801
    43 = 2
    45 = 38
    46 = 484
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
